$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.405.00'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").Value = '1.905.27'
$ws.Range("E3").Value = '  -2.85%  '
$ws.Range("D4").Value = '''0.9986'
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").Value = '''239.23'
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("D6").Value = '''0.9989'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").Value = '''0.4728'
$ws.Range("E7").Value = '  -2.50%  '
$ws.Range("D8").Value = '''0.2835'
$ws.Range("E8").Value = '  -3.82%  '
$ws.Range("D9").Value = '''0.06660'
$ws.Range("E9").Value = '  -5.85%  '
$ws.Range("D10").Value = '''18.71'
$ws.Range("E10").Value = '  -5.00%  '
$ws.Range("D11").Value = '''99.64'
$ws.Range("E11").Value = '  -6.98%  '
$ws.Range("D12").Value = '''0.07713'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '1.913.99'
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("D14").Value = '''5.200'
$ws.Range("E14").Value = '  -3.76%  '
$ws.Range("D15").Value = '''0.6683'
$ws.Range("E15").Value = '  -5.34%  '
$ws.Range("D16").Value = '30.425.29'
$ws.Range("E16").Value = '  -2.04%  '
$ws.Range("D17").Value = '''254.24'
$ws.Range("E17").Value = '  -8.90%  '
$ws.Range("D18").Value = '''0.9989'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").Value = '''0.000007438'
$ws.Range("E19").Value = '  -4.95%  '
$ws.Range("E20").Value = '  -5.58%  '
$ws.Range("D21").Value = '''5.371'
$ws.Range("E21").Value = '  -2.96%  '
$ws.Range("D22").Value = '''0.9984'
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").Value = '''0.4504'
$ws.Range("E23").Value = '  -8.45%  '
$ws.Range("D24").Value = '''6.314'
$ws.Range("E24").Value = '  -3.17%  '
$ws.Range("D25").Value = '''9.434'
$ws.Range("E25").Value = '  -3.55%  '
$ws.Range("D26").Value = '''167.37'
$ws.Range("E26").Value = '  -1.17%  '
$ws.Range("D27").Value = '''18.86'
$ws.Range("E27").Value = '  -4.57%  '
$ws.Range("D28").Value = '''2.048'
$ws.Range("E28").Value = '  -6.21%  '
$ws.Range("E29").Value = '  -4.09%  '
$ws.Range("D30").Value = '''1.375'
$ws.Range("E30").Value = '  -1.91%  '
$ws.Range("D31").Value = '''4.630'
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  -3.71%  '
$ws.Range("D33").Value = '''4.250'
$ws.Range("E33").Value = '  -4.11%  '
$ws.Range("D34").Value = '''0.04712'
$ws.Range("E34").Value = '  -4.07%  '
$ws.Range("D35").Value = '''0.7284'
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("E36").Value = '  -5.50%  '
$ws.Range("D37").Value = '''0.9981'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").Value = '''2.696'
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").Value = '''0.01915'
$ws.Range("E39").Value = '  -4.81%  '
$ws.Range("D40").Value = '''2.595'
$ws.Range("E40").Value = '  -3.50%  '
$ws.Range("E41").Value = '  -4.01%  '
$ws.Range("D42").Value = '''73.54'
$ws.Range("E42").Value = '  -6.89%  '
$ws.Range("D43").Value = '''1.958'
$ws.Range("E43").Value = '  -8.34%  '
$ws.Range("D44").Value = '''0.8586'
$ws.Range("E44").Value = '  -4.42%  '
$ws.Range("D45").Value = '''105.55'
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("D46").Value = '''0.9982'
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").Value = '''0.4225'
$ws.Range("E47").Value = '  -5.22%  '
$ws.Range("D48").Value = '''7.393'
$ws.Range("E48").Value = '  -6.70%  '
$ws.Range("D49").Value = '''977.31'
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").Value = '''0.1196'
$ws.Range("E50").Value = '  -4.44%  '
$ws.Range("D51").Value = '''34.58'
$ws.Range("E51").Value = '  -4.08%  '
